# Update gh-pages to output generated at 456a3b4
# Applies numeric "想去人数" (interest count) bumps across the 展览 / 演出 /
# 全部类型 sheets, relocates the 李宁运动中心 venue text, and replaces the
# last exhibition row (奶司的小人国娃展) with a new 第五人格ONLY entry while
# pushing the old row down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$sheet1Changes = @{
    2  = 929
    3  = 1029
    4  = 819
    7  = 725
    8  = 171
    9  = 1324
    10 = 752
    12 = 569
    13 = 192
    14 = 60
    15 = 60
    16 = 1254
    17 = 150
    20 = 381
    21 = 98
    22 = 608
    23 = 167
    24 = 667
    25 = 40
    26 = 1113
}
foreach ($row in $sheet1Changes.Keys) {
    $ws1.Cells.Item([int]$row, 6).Value = $sheet1Changes[$row]
}

$ws1.Cells.Item(19, 4).Value = "广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心"

# Row 27 (奶司的小人国娃展) becomes row 28, and a brand-new row 27
# (第五人格ONLY) is inserted in its place.
$ws1.Rows("27:27").Insert()
$ws1.Cells.Item(27, 1).Value = 26
$ws1.Cells.Item(27, 1).Borders.LineStyle = 1

$ws1.Cells.Item(27, 2).NumberFormat = "@"
$ws1.Cells.Item(27, 2).Value = "2024-05-18"
$ws1.Cells.Item(27, 3).Value = "广州·第五人格ONLY"
$ws1.Cells.Item(27, 4).Value = "洛浦街厦滘西环路1号 广州市岭南国际电子商务会展中心"
$ws1.Cells.Item(27, 5).Value = "2024.05.18 10:00-05.18 17:00"
$ws1.Cells.Item(27, 6).Value = 0
$ws1.Cells.Item(27, 7).Value = 60
$ws1.Cells.Item(27, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82458"
$ws1.Cells.Item(27, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/D8jK0O2X1709778592031.jpeg"

# The row that got pushed down (now row 28) keeps its own sequential id,
# bumped by one to 27.
$ws1.Cells.Item(28, 1).Value = 27

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$sheet2Changes = @{
    8  = 59
    11 = 120
}
foreach ($row in $sheet2Changes.Keys) {
    $ws2.Cells.Item([int]$row, 6).Value = $sheet2Changes[$row]
}

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) -- mirrors sheet 1 with a wider row range
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$sheet4Changes = @{
    4  = 929
    5  = 1029
    6  = 819
    9  = 725
    10 = 171
    11 = 1324
    12 = 752
    16 = 569
    18 = 192
    19 = 60
    20 = 60
    21 = 1254
    23 = 150
    26 = 381
    27 = 98
    29 = 59
    30 = 608
    33 = 120
    34 = 120
    35 = 167
    36 = 667
    37 = 40
    38 = 1113
}
foreach ($row in $sheet4Changes.Keys) {
    $ws4.Cells.Item([int]$row, 6).Value = $sheet4Changes[$row]
}

$ws4.Cells.Item(25, 4).Value = "广龙路中油BP(白云万顺达南加油站)北侧约260米 李宁运动中心"

# Row 39 (奶司的小人国娃展) becomes row 40, and a brand-new row 39
# (第五人格ONLY) is inserted in its place.
$ws4.Rows("39:39").Insert()
$ws4.Cells.Item(39, 1).Value = 38
$ws4.Cells.Item(39, 1).Borders.LineStyle = 1

$ws4.Cells.Item(39, 2).NumberFormat = "@"
$ws4.Cells.Item(39, 2).Value = "2024-05-18"
$ws4.Cells.Item(39, 3).Value = "广州·第五人格ONLY"
$ws4.Cells.Item(39, 4).Value = "洛浦街厦滘西环路1号 广州市岭南国际电子商务会展中心"
$ws4.Cells.Item(39, 5).Value = "2024.05.18 10:00-05.18 17:00"
$ws4.Cells.Item(39, 6).Value = 0
$ws4.Cells.Item(39, 7).Value = 60
$ws4.Cells.Item(39, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82458"
$ws4.Cells.Item(39, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/D8jK0O2X1709778592031.jpeg"

# The row that got pushed down (now row 40) keeps its own sequential id,
# bumped by one to 39.
$ws4.Cells.Item(40, 1).Value = 39

Write-Host "edit.ps1 applied"
